# Generate Report for Handoff
# The CI job re-ran for a renamed source file (new GUID) and refreshed the
# handoff timestamps. Because the handback step has not happened yet for
# the new file, the "latest target / handback file / handback datetime"
# columns for each locale sheet are reset.

$wb = $excel.ActiveWorkbook

$oldGuid = "e5b5e060-5285-4db9-8e1f-6904a89979b8"
$newGuid = "49f12dad-2df1-4f2d-b999-df0a7667fe4e"
$newHash = "e734198ed52e1edfcf2e67168d5e7de5ecf22772"

$oldFileName = "$oldGuid.md"
$newFileName = "$newGuid.md"
$oldPathName = "e2e\$oldGuid.md"
$newPathName = "e2e\$newGuid.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value = $newFileName
$overview.Range("B2").Value = $newPathName
$overview.Range("G2").Value = "2016-09-04 11:05:38"

$overviewUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6dcaff99921fb9d2b2d01b3143d5bde785797242/e2e/$oldGuid.md"
$overview.Range("B2").Hyperlinks.Delete()
$overview.Hyperlinks.Add($overview.Range("B2"), $overviewUrl, [System.Type]::Missing, [System.Type]::Missing, $newPathName) | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value = $newFileName
$zhcn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-09-04 11:05:33"
$zhcn.Range("I2").Value = ""
$zhcn.Range("I2").Style = "Normal"
$zhcn.Range("J2").Value = ""
$zhcn.Range("J2").Style = "Normal"
$zhcn.Range("K2").Value = "0001-01-01 00:00:00"

$zhcnUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6dcaff99921fb9d2b2d01b3143d5bde785797242/e2e/$oldGuid.md"
$zhcn.Range("A2").Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $zhcnUrl, [System.Type]::Missing, [System.Type]::Missing, $newFileName) | Out-Null

$zhcn.Columns.Item(9).ColumnWidth = 18.6506053379604
$zhcn.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value = $newFileName
$dede.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$dede.Range("H2").Value = "2016-09-04 11:05:38"
$dede.Range("I2").Value = ""
$dede.Range("I2").Style = "Normal"
$dede.Range("J2").Value = ""
$dede.Range("J2").Style = "Normal"
$dede.Range("K2").Value = "0001-01-01 00:00:00"

$dedeUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6dcaff99921fb9d2b2d01b3143d5bde785797242/e2e/$oldGuid.md"
$dede.Range("A2").Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $dedeUrl, [System.Type]::Missing, [System.Type]::Missing, $newFileName) | Out-Null

$dede.Columns.Item(9).ColumnWidth = 18.6506053379604
$dede.Columns.Item(10).ColumnWidth = 21.7054770333426
